# Spotify/Code.xlsx update: add a new "405" row (Tham gia Premium Family
# tu mot quoc gia khac) right above the existing "408" row, pushing the
# 408/409 rows down by one, widen column D, and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 9 (the "408" row); rows 9-10 become 10-11.
$ws.Rows("9:9").Insert()

# Give the new row the same look (borders + centered text) as the rest of
# the table by copying the formatting from the row right above it.
$ws.Range("B8:D8").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new row's content.
$ws.Range("B9").Value = 405
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "Tham gia Premium Family từ một quốc gia khác"

# Widen column D (37 -> 46.5 characters).
$ws.Columns("D").ColumnWidth = 45.785714285714285

# Match the saved selection/active cell.
$ws.Range("J11").Select() | Out-Null
